$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 5 ---
$ws.Range("C5").Value = "13:00"
$ws.Range("D5").Value = "15:00"
$ws.Range("E5").Value = "aaa"
$ws.Range("F5").Value = "Sandy"
$ws.Range("G5").Value = "Honeywell"
$ws.Range("H5").Value = "CCC"

# --- New row 6 ---
$ws.Range("A6").Value = "sandy"
$ws.Range("B6").Value = "'2025-03-27"
$ws.Range("C6").Value = "14:00"
$ws.Range("D6").Value = "15:00"
$ws.Range("E6").Value = "vellore"
$ws.Range("F6").Value = "Sandy"
$ws.Range("G6").Value = "Honeywell"
$ws.Range("H6").Value = "CCC"

# --- New row 7 ---
$ws.Range("A7").Value = "Mandy"
$ws.Range("B7").Value = "'2024-03-26"
$ws.Range("C7").Value = "16:01"
$ws.Range("D7").Value = "17:00"
$ws.Range("E7").Value = "test"
$ws.Range("F7").Value = "test"
$ws.Range("G7").Value = "test"
$ws.Range("H7").Value = "test"
